$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7

$ws.Cells.Item($row, 1).Value = 8
$ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 44964
$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = 100114002
$ws.Cells.Item($row, 7).Value = "Camote"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 300
$ws.Cells.Item($row, 11).Value = 20000
$ws.Cells.Item($row, 12).Value = 21000
$ws.Cells.Item($row, 13).Value = 20500
$ws.Cells.Item($row, 14).Value = '$/malla 18 kilos'
$ws.Cells.Item($row, 15).Value = "Perú"
$ws.Cells.Item($row, 16).Value = 1139
$ws.Cells.Item($row, 17).Value = 18
$ws.Cells.Item($row, 18).Value = "Hortaliza"

$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
